$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data correction ---
# C2 and C6 were miscoded ("A" / "B") and should read "C", matching the
# grading values already used on several other rows.
$ws.Range("C2").Value = "C"
$ws.Rows(2).AutoFit()
$ws.Range("C6").Value = "C"
$ws.Rows(6).AutoFit()

# --- Refresh / move the AutoFilter so it covers the newly added row 43 ---
# Clear the existing filter first so Excel recalculates the range from the
# sheet's used range (A1:C43) instead of keeping the stale A1:C42 ref.
$ws.AutoFilterMode = $false

# Re-apply the AutoFilter over the full data range, now filtering on
# Program values 1 and 2 (instead of 37).
$ur = $ws.UsedRange
$ur.AutoFilter(1, @("1","2"), 7)

# Keep the workbook-level _FilterDatabase defined name in sync with the
# new AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$43"
    }
}

# --- Update the last active selection on the sheet ---
$ws.Range("C3").Select()
